$wb = $excel.ActiveWorkbook

# Sheet1: update distributor value in B8 (new shared string added first so it
# lands at the lower index) and move the selection to A10.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("B8").Value = "Fosroc@3"
$ws1.Range("A10").Select()

# SEBS_Devloper: add the new distributor name in A3 and move the selection to
# C6. Re-activate Sheet1 afterwards so it remains the selected/active tab.
$ws2 = $wb.Worksheets.Item("SEBS_Devloper")
$ws2.Activate()
$ws2.Range("A3").Value = "sahil Khan"
$ws2.Range("C6").Select()

$ws1.Activate()
$ws1.Range("A10").Select()
